# -------------------------------------------------------------------------
# contratos-2-2009.xlsx formatting fix
#
# 1) Four "Razon social"/"Nombre Fantasia" text entries used "," as a
#    separator between co-titulares; the scraper that produced this sheet
#    mangled them (comma -> period) the same way it mangled decimals below.
# 2) The "Importe" column (H) was scraped as es-AR formatted text
#    ("1.234,56" = thousands "." + decimal ",") but needed to read as plain
#    "1234.56" (decimal "." , no thousands separator).
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "Razon social" / "Nombre Fantasia" punctuation fixes ---------------

$ws.Range("E32").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F32").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E62").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F62").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'

$ws.Range("E53").Value = 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'
$ws.Range("F53").Value = 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'

$ws.Range("E61").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'

$ws.Range("E63").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'

# --- 2) "Importe" column (H): es-AR "1.234,56" -> plain "1234.56" ----------
#     Column H cells are stored as TEXT (not numbers) in the source sheet,
#     so each new value is re-entered as text (NumberFormat "@" beforehand
#     keeps Excel from re-parsing the digit string back into a real number;
#     resetting the style back to Normal afterwards keeps the cell format
#     identical to how it started).

$importeFixes = @{
    2 = "340.00"
    3 = "555.00"
    4 = "79.00"
    5 = "2715.00"
    6 = "187.40"
    7 = "736.00"
    8 = "54.00"
    9 = "15964.38"
    10 = "46244.62"
    11 = "8737.71"
    12 = "5557.55"
    13 = "10453.87"
    14 = "62.50"
    15 = "4451.39"
    16 = "3541.30"
    17 = "5262.37"
    18 = "216.63"
    19 = "6427.83"
    20 = "4088.50"
    21 = "3270.00"
    22 = "61.05"
    23 = "1490.00"
    24 = "241.00"
    25 = "58.34"
    26 = "1807.13"
    27 = "3493.95"
    28 = "4858.35"
    29 = "42039.00"
    30 = "5973.00"
    31 = "86.15"
    32 = "3919.65"
    33 = "19965.20"
    34 = "2403.92"
    35 = "264.00"
    36 = "12937.20"
    37 = "8.32"
    38 = "15300.00"
    39 = "1841.22"
    40 = "270.00"
    41 = "5013.70"
    42 = "369.00"
    43 = "7931.86"
    44 = "3730.00"
    45 = "894.39"
    46 = "716.00"
    47 = "9873.60"
    48 = "800.00"
    49 = "1600.00"
    50 = "1003.40"
    51 = "8015.00"
    52 = "800.00"
    53 = "215.00"
    54 = "807.50"
    55 = "12848.40"
    56 = "14875.00"
    57 = "2000.00"
    58 = "509.00"
    59 = "32941.00"
    60 = "1706.60"
    61 = "108.00"
    62 = "442.57"
    63 = "2071.00"
    64 = "680.00"
    65 = "15.00"
    66 = "213532.00"
    67 = "262.29"
    68 = "239.61"
    69 = "3334.00"
    70 = "154.17"
    71 = "433.07"
    72 = "270.00"
    73 = "388.50"
    74 = "2455.52"
    75 = "5807.00"
    76 = "4484.34"
    77 = "308.50"
    78 = "2166.80"
    79 = "6249.39"
    80 = "210.00"
    81 = "46.89"
    82 = "6960.13"
    83 = "24581.46"
    84 = "354.00"
    85 = "964.00"
    86 = "290.00"
    87 = "2100.00"
    88 = "1500.00"
    89 = "1040.00"
    90 = "233.00"
    91 = "2450.00"
    92 = "600.00"
    93 = "500.00"
    94 = "830.00"
    95 = "597.18"
    96 = "62051.22"
    97 = "574.60"
    98 = "400.00"
    99 = "500.00"
    100 = "700.00"
    101 = "150.00"
    102 = "756.00"
    103 = "433.16"
    104 = "3025.00"
    105 = "240.00"
    106 = "250.00"
    107 = "728.00"
    108 = "4869.74"
    109 = "1350.00"
    110 = "1000.00"
    111 = "150.00"
    112 = "720.00"
    113 = "13699.62"
    114 = "700.00"
    115 = "494.00"
    116 = "240.00"
    117 = "3000.00"
    118 = "1456.85"
    119 = "362.00"
    120 = "2977.50"
    121 = "108.90"
    122 = "1354.56"
    123 = "2182.45"
    124 = "245.00"
    125 = "3200.00"
    126 = "2200.00"
    127 = "229.84"
    128 = "2430.78"
    129 = "144.00"
    130 = "298.63"
    131 = "261.00"
    132 = "127.80"
    133 = "372.00"
    134 = "12641.30"
    135 = "599.00"
    136 = "195.00"
    137 = "5801.96"
    138 = "3280.00"
    139 = "1248.96"
    140 = "2120.00"
    141 = "987.00"
    142 = "2001.52"
    143 = "4125.00"
    144 = "197482.31"
    145 = "1552.30"
    146 = "490.00"
    147 = "2240.00"
}

foreach ($row in $importeFixes.Keys) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.NumberFormat = "@"
    $cell.Value = $importeFixes[$row]
    $cell.Style = "Normal"
}

